$wb = $excel.ActiveWorkbook

# --- Update the "Date" value on the Metadata sheet ---
$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B8").Value = "2023-02-21T12:43:59+00:00"

# --- Update the EXOR -> REPORT concept row on the Concepts sheet ---
$concepts = $wb.Worksheets.Item("Concepts")
$concepts.Range("B8").Value = "REPORT"
$concepts.Range("C8").Value = "Data Report"
